# Scheduled market-data refresh for Rafflesia_Profits leve-profit sheets.
# Refreshes the cached currentAveragePrice / currentAveragePriceNQ /
# currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ /
# LeveProfitHQ columns (H:N) for the leves whose market snapshot changed
# since the previous run. Rows/leves not listed below are left untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9 (Leve Item ID 5487)
$ws.Range("H9").Value = 151.57143
$ws.Range("I9").Value = 157.6
$ws.Range("J9").Value = 136.5
$ws.Range("K9").Value = 157.6
$ws.Range("L9").Value = 136.5
$ws.Range("M9").Value = 11.40000000000001
$ws.Range("N9").Value = -474.5

# Row 17 (Leve Item ID 38956)
$ws.Range("H17").Value = 2785
$ws.Range("J17").Value = 2785
$ws.Range("L17").Value = 8355
$ws.Range("N17").Value = -8691

# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 2444
$ws.Range("I43").Value = 2444
$ws.Range("K43").Value = 2444
$ws.Range("M43").Value = -2375

# Row 51 (Leve Item ID 5486)
$ws.Range("H51").Value = 10055.546
$ws.Range("I51").Value = 10666.667
$ws.Range("J51").Value = 9826.375
$ws.Range("K51").Value = 10666.667
$ws.Range("L51").Value = 9826.375
$ws.Range("M51").Value = -10182.667
$ws.Range("N51").Value = -10794.375

# Row 92 (Leve Item ID 19901)
$ws.Range("H92").Value = 222.8
$ws.Range("I92").Value = 277.25
$ws.Range("J92").Value = 5
$ws.Range("K92").Value = 277.25
$ws.Range("L92").Value = 5
$ws.Range("M92").Value = 970.75
$ws.Range("N92").Value = -2501

# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 747.4
$ws.Range("I98").Value = 412.66666
$ws.Range("J98").Value = 1249.5
$ws.Range("K98").Value = 412.66666
$ws.Range("L98").Value = 1249.5
$ws.Range("M98").Value = 1085.33334
$ws.Range("N98").Value = -4245.5

# Row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 378.33334
$ws.Range("I107").Value = 70
$ws.Range("J107").Value = 686.6667
$ws.Range("K107").Value = 70
$ws.Range("L107").Value = 686.6667
$ws.Range("M107").Value = 1850
$ws.Range("N107").Value = -4526.6667

# Row 113 (Leve Item ID 27775)
$ws.Range("H113").Value = 3561.25
$ws.Range("I113").Value = 3561.25
$ws.Range("K113").Value = 3561.25
$ws.Range("M113").Value = -307.25

# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 747.4
$ws.Range("I122").Value = 412.66666
$ws.Range("J122").Value = 1249.5
$ws.Range("K122").Value = 1237.99998
$ws.Range("L122").Value = 3748.5
$ws.Range("M122").Value = 1212.00002
$ws.Range("N122").Value = -8648.5

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1386.2858
$ws.Range("I137").Value = 1317.3334
$ws.Range("K137").Value = 3952.0002
$ws.Range("M137").Value = -1402.0002

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2150
$ws.Range("J138").Value = 2850
$ws.Range("L138").Value = 8550
$ws.Range("N138").Value = -18830

$ws = $wb.Worksheets.Item("ARM")
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 4005
$ws.Range("J132").Value = 10000
$ws.Range("L132").Value = 30000
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("BSM")
# Row 100 (Leve Item ID 18347)
$ws.Range("H100").Value = 20643
$ws.Range("J100").Value = 20643
$ws.Range("L100").Value = 20643
$ws.Range("N100").Value = -22807

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 6250
$ws.Range("I134").Value = 1562.5
$ws.Range("K134").Value = 4687.5
$ws.Range("M134").Value = -2152.5

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

# Row 28 (Leve Item ID 18348)
$ws.Range("H28").Value = 42504.5
$ws.Range("J28").Value = 42504.5
$ws.Range("L28").Value = 42504.5
$ws.Range("N28").Value = -42994.5

$ws = $wb.Worksheets.Item("CUL")
# Row 12 (Leve Item ID 4854)
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 719.2857
$ws.Range("I107").Value = 267.66666
$ws.Range("J107").Value = 1058
$ws.Range("K107").Value = 802.9999799999999
$ws.Range("L107").Value = 3174
$ws.Range("M107").Value = 1117.00002
$ws.Range("N107").Value = -7014

# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 2282.353
$ws.Range("J131").Value = 2312.5
$ws.Range("L131").Value = 6937.5
$ws.Range("N131").Value = -17017.5

# Row 137 (Leve Item ID 44088)
$ws.Range("H137").Value = 20000
$ws.Range("J137").Value = 20000
$ws.Range("L137").Value = 60000
$ws.Range("N137").Value = -70200

# Row 138 (Leve Item ID 44105)
$ws.Range("H138").Value = 1680.625
$ws.Range("J138").Value = 2000
$ws.Range("L138").Value = 6000
$ws.Range("N138").Value = -16280

$ws = $wb.Worksheets.Item("GSM")
# Row 95 (Leve Item ID 18235)
$ws.Range("H95").Value = 27374.5
$ws.Range("J95").Value = 27374.5
$ws.Range("L95").Value = 27374.5
$ws.Range("N95").Value = -32866.5

# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 1262
$ws.Range("I122").Value = 1262
$ws.Range("K122").Value = 3786
$ws.Range("M122").Value = -1336

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 6506
$ws.Range("I126").Value = 6506
$ws.Range("K126").Value = 19518
$ws.Range("M126").Value = -17048

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 2675
$ws.Range("I132").Value = 2012.5
$ws.Range("K132").Value = 6037.5
$ws.Range("M132").Value = -3507.5

$ws = $wb.Worksheets.Item("LTW")
# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 1200
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1200
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 1200
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -2698

# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 1200
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1200
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 6000
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -13488

# Row 94 (Leve Item ID 18067)
$ws.Range("H94").Value = 53661.332
$ws.Range("J94").Value = 53661.332
$ws.Range("L94").Value = 53661.332
$ws.Range("N94").Value = -55013.332

# Row 95 (Leve Item ID 18221)
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# Row 104 (Leve Item ID 18675)
$ws.Range("H104").Value = 32038.75
$ws.Range("J104").Value = 32038.75
$ws.Range("L104").Value = 32038.75
$ws.Range("N104").Value = -39026.75

$ws = $wb.Worksheets.Item("WVR")
# Row 45 (Leve Item ID 21726)
$ws.Range("H45").Value = 30626
$ws.Range("J45").Value = 30626
$ws.Range("L45").Value = 30626
$ws.Range("N45").Value = -31608

# Row 62 (Leve Item ID 12589)
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# Row 65 (Leve Item ID 12589)
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 555
$ws.Range("J107").Value = 555
$ws.Range("L107").Value = 1665
$ws.Range("N107").Value = -5505
